$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row formatting: bold font, thin box border, centered + top-aligned ---
$header = $ws.Range("A1:K1")
$header.Font.Bold = $true
$header.Borders.LineStyle = 1
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4160

# --- Row 2 updates ---
$ws.Range("C2").Value = 2
$ws.Range("D2").Value = "2026-02-12T10:14:26.458480+00:00"
$ws.Range("E2").Value = "Yes, could you share their contacts"
$ws.Range("F2").Value = "yes"

# --- Row 3 updates ---
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = "2026-02-12T10:14:54.857564+00:00"
$ws.Range("E3").Value = "Yes I am looking for opportunities"
$ws.Range("F3").Value = "yes"

Write-Host "Applied header style and updated rows 2-3"
